# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for the rows
# whose dialog-act annotation changed after re-running SGNN on the cleaned transcripts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    22  = @("aa", "Agree/Accept")
    29  = @("sv", "Statement-opinion")
    46  = @("sd", "Statement-non-opinion")
    47  = @("b", "Acknowledge (Backchannel)")
    55  = @("sd", "Statement-non-opinion")
    61  = @("aa", "Agree/Accept")
    113 = @("aa", "Agree/Accept")
    118 = @("ba", "Appreciation")
    119 = @("sd", "Statement-non-opinion")
    145 = @("sd", "Statement-non-opinion")
    146 = @("sd", "Statement-non-opinion")
    152 = @("%", "Uninterpretable")
    159 = @("sd", "Statement-non-opinion")
    166 = @("sd", "Statement-non-opinion")
    172 = @("b", "Acknowledge (Backchannel)")
    176 = @("b", "Acknowledge (Backchannel)")
    191 = @("sv", "Statement-opinion")
    199 = @("sv", "Statement-opinion")
    200 = @("b", "Acknowledge (Backchannel)")
    201 = @("sv", "Statement-opinion")
    203 = @("ba", "Appreciation")
    215 = @("b", "Acknowledge (Backchannel)")
    253 = @("sd", "Statement-non-opinion")
    255 = @("sd", "Statement-non-opinion")
    259 = @("%", "Uninterpretable")
    269 = @("%", "Uninterpretable")
    296 = @("sd", "Statement-non-opinion")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Range("I$row").Value = $values[0]
    $ws.Range("J$row").Value = $values[1]
}
